# Generate Report for Handback
#
# This script mutates the localization-status workbook the way the
# OpenLocalization "handback" job does once a target (de-de) translation
# has round-tripped and is back in sync with en-US:
#   * the "Ready for handoff" status becomes "Handed back: in sync with en-US"
#     everywhere it is used (Overview + per-language sheets),
#   * each language sheet's row gets its "Latest Target File" (I) and
#     "Latest Handback File" (J) columns populated, and "Latest Handback
#     DateTime" (K) stamped,
#   * the Source-Path/Target-File columns are widened to fit the longer
#     file names now appearing in them.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddd520aa1a543a0628de92e6165f6f122a0efa04/e2e/"
$mdA = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.md"
$mdB = "ffff778366a4-ef15-40b7-a639-bcc0e5053014.md"

$newStatus = "Handed back: in sync with en-US"

$ovw   = $wb.Worksheets.Item("Overview")
$zhcn  = $wb.Worksheets.Item("zh-cn")
$dede  = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    on every cell that currently shows it.
# ---------------------------------------------------------------------
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: record the target/handback file + handback datetime.
#    "Latest Target File" and "Latest Handback File" both link back to
#    the source markdown file (same as column A), and the handback
#    timestamp is stamped in column K.
# ---------------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($baseUrl + $mdA), [Type]::Missing, [Type]::Missing, $mdA)
$zhcn.Range("J2").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 07:14:20"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($baseUrl + $mdA), [Type]::Missing, [Type]::Missing, $mdA)
$zhcn.Range("J3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 07:14:20"

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape, with the de-de xlf + its own handback time.
# ---------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), ($baseUrl + $mdA), [Type]::Missing, [Type]::Missing, $mdA)
$dede.Range("J2").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 07:14:28"

$dede.Hyperlinks.Add($dede.Range("I3"), ($baseUrl + $mdA), [Type]::Missing, [Type]::Missing, $mdA)
$dede.Range("J3").Value = "a9da29a6-e88c-4baf-9ba0-43eaec9ccccb.4c0a76b601f5a3a4ae2cd76245157bf4ac510dbf.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 07:14:28"

# ---------------------------------------------------------------------
# 4. Widen the columns that now hold the longer file-name / datetime
#    strings, on all three sheets.
# ---------------------------------------------------------------------
$ovw.Range("E1").ColumnWidth = 29.9777047293527
$ovw.Range("F1").ColumnWidth = 29.9777047293527

$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40
